$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A previously spanned cols 1-2 (shadowed by a later, more specific rule for
# col B); narrow its span down to col 1 only. Toggling Hidden forces the engine to split
# the merged column-format span without perturbing the stored width value.
$ws.Columns("A:A").Hidden = $true
$ws.Columns("A:A").Hidden = $false

# --- Objetivos (row 10): replace the placeholder value with the full objectives paragraph
$ws.Range("B10").Value = "Introduzir o estudante em conceitos importantes utilizados nas ciências microbiológicas, apresentando os conceitos fundamentais relacionados à história, mercado, genética, filogenia, e cultivo de microrganismos."
$ws.Range("C10").Value = "Introduzir o estudante em conceitos importantes utilizados nas ciências microbiológicas, apresentando os conceitos fundamentais relacionados à história, mercado, genética, filogenia, e cultivo de microrganismos."

# --- Insert two rows under "Docentes responsaveis:" to list both professors
$ws.Rows("13:14").Insert()
$ws.Range("B13").Value = "4873328 - Fernando Segato"
$ws.Range("C13").Value = "4873328 - Fernando Segato"
$ws.Range("B14").Value = "5840685 - Maria Bernadete de Medeiros"
$ws.Range("C14").Value = "5840685 - Maria Bernadete de Medeiros"

# The inserted rows picked up column A formatting from the row above; this section has no A value
$ws.Range("A13:A14").Clear()

# Re-apply the correct B/C formatting (wrap-top text / wrap-top red text) to the new rows
$ws.Range("B16").Copy()
$ws.Range("B13:B14").PasteSpecial(-4122)
$ws.Range("C16").Copy()
$ws.Range("C13:C14").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Programa resumido (row 15, was 13): fix the short-syllabus PT text
$ws.Range("B15").Value = "Histórico da microbiologia, microbiologia industrial, filogênia microbiana, caracterização dos microrganismos, nutrição e cultivo de microrganismos, virus, fungos filamentosos, leveduras, micro-algas, bactérias."
$ws.Range("C15").Value = "Histórico da microbiologia, microbiologia industrial, filogênia microbiana, caracterização dos microrganismos, nutrição e cultivo de microrganismos, virus, fungos filamentosos, leveduras, micro-algas, bactérias."

# --- Programa (row 17, was 15): fix the long PT program text
$ws.Range("B17").Value = "1. Histórico da microbiologia: microbiologia, ciência e sociedade;  Leeuwenhoek e seusseus microscópios, origem dos animálculos de Leeuwenhoek, principais pensadores da microbiologia, microbiologia moderna.2. Microbiologia industrial: visão geral do mercado envolvendo microbiologia, principais produtos de origem microbiana.3. Filogênia microbiana: classificação e evolução das principais classes dos microrganismos; organismos procarióticos (eubactérias e arqueobactérias); organismos eucarióticos (leveduras, fungos filamentosos, algas, protozoários).4. Caracterização dos microrganismos: técnicas de cultura pura, microscópios, técnicas de microscopia, preparo dos microrganismos para microscopia, informações utilizadas para caracterizar os microrganismos, tecnologia automatizada.5. Nutrição e cultivo de microrganismos: exigências nutricionais e meios microbiológicos; cultivo ecrescimento de microrganismos.6. Genética de microrganismos: regulação da expressão gênica em bactérias; mutação, vantagens e desvantagens para aplicações industriais; melhoramento de cepas.7. Virus, bactérias, fungos filamentosos, micro-algas, leveduras: morfologia, classificação e replicação."
$ws.Range("C17").Value = "1. Histórico da microbiologia: microbiologia, ciência e sociedade;  Leeuwenhoek e seusseus microscópios, origem dos animálculos de Leeuwenhoek, principais pensadores da microbiologia, microbiologia moderna.2. Microbiologia industrial: visão geral do mercado envolvendo microbiologia, principais produtos de origem microbiana.3. Filogênia microbiana: classificação e evolução das principais classes dos microrganismos; organismos procarióticos (eubactérias e arqueobactérias); organismos eucarióticos (leveduras, fungos filamentosos, algas, protozoários).4. Caracterização dos microrganismos: técnicas de cultura pura, microscópios, técnicas de microscopia, preparo dos microrganismos para microscopia, informações utilizadas para caracterizar os microrganismos, tecnologia automatizada.5. Nutrição e cultivo de microrganismos: exigências nutricionais e meios microbiológicos; cultivo ecrescimento de microrganismos.6. Genética de microrganismos: regulação da expressão gênica em bactérias; mutação, vantagens e desvantagens para aplicações industriais; melhoramento de cepas.7. Virus, bactérias, fungos filamentosos, micro-algas, leveduras: morfologia, classificação e replicação."

# --- Metodo (row 20, was 18): fix the evaluation-method text
$ws.Range("B20").Value = "A avaliação será feita por meio de provas escritas, trabalhos, seminários e participação."
$ws.Range("C20").Value = "A avaliação será feita por meio de provas escritas, trabalhos, seminários e participação."

# --- Criterio (row 21, was 19): fix the final-grade formula text
$ws.Range("B21").Value = "A Nota final (NF) será calculada da seguinte maneira: NF = (P1 + P2)/2."
$ws.Range("C21").Value = "A Nota final (NF) será calculada da seguinte maneira: NF = (P1 + P2)/2."

# --- Norma de recuperacao (row 22, was 20): fix the recovery-norm text
$ws.Range("B22").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR)"
$ws.Range("C22").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR)"

# --- Bibliografia (row 23, was 21): fix the bibliography text
$ws.Range("B23").Value = "1. PELCZAR Jr, M.J., CHAN, S.S., KRIEG, N.R. Microbiologia conceitos e aplicações, 2 ed. (Vol 1), São Paulo: Pearson Education do Brasil, 1997.2. MADIGAN, M.T., MARTINKO, J.M., PARKER, I. Microbiologia de Brock. São Paulo: Prentice Hall, 2004.3. BARBOSA, H.R., TORRES, B.B. Microbiologia Básica, Rio de Janeiro: Atheneu, 2005.4. VERMELHO A.B., FREIRE, M.C., BRANQUINHO, M.H. Bacteorologia Geral, Rio de Janeiro: GuanabaraKoogan, 2008.5. TORTORA, G.J., FUNKE, B.R., CASE, C.L. Microbiologia, Artmed, Porto Alegre, RS, 2012."
$ws.Range("C23").Value = "1. PELCZAR Jr, M.J., CHAN, S.S., KRIEG, N.R. Microbiologia conceitos e aplicações, 2 ed. (Vol 1), São Paulo: Pearson Education do Brasil, 1997.2. MADIGAN, M.T., MARTINKO, J.M., PARKER, I. Microbiologia de Brock. São Paulo: Prentice Hall, 2004.3. BARBOSA, H.R., TORRES, B.B. Microbiologia Básica, Rio de Janeiro: Atheneu, 2005.4. VERMELHO A.B., FREIRE, M.C., BRANQUINHO, M.H. Bacteorologia Geral, Rio de Janeiro: GuanabaraKoogan, 2008.5. TORTORA, G.J., FUNKE, B.R., CASE, C.L. Microbiologia, Artmed, Porto Alegre, RS, 2012."

